$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-07 Sunday" "2025-09-08 Monday"

Replace-Text "431÷2=215, 1" "487÷9=54, 1"
Replace-Text "697÷3=232, 1" "693÷6=115, 3"
Replace-Text "357÷2=178, 1" "149÷4=37, 1"
Replace-Text "737÷3=245, 2" "131÷6=21, 5"
Replace-Text "286÷9=31, 7" "176÷3=58, 2"

Replace-Text "436÷5=87, 1" "483÷3=161, 0"
Replace-Text "541÷4=135, 1" "862÷7=123, 1"
Replace-Text "378÷5=75, 3" "401÷3=133, 2"
Replace-Text "713÷7=101, 6" "603÷9=67, 0"
Replace-Text "472÷7=67, 3" "735÷2=367, 1"

Replace-Text "191÷8=23, 7" "669÷2=334, 1"
Replace-Text "520÷8=65, 0" "555÷7=79, 2"
Replace-Text "924÷7=132, 0" "581÷5=116, 1"
Replace-Text "504÷5=100, 4" "160÷9=17, 7"
Replace-Text "291÷5=58, 1" "658÷3=219, 1"

Replace-Text "370÷8=46, 2" "152÷4=38, 0"
Replace-Text "613÷7=87, 4" "992÷9=110, 2"
Replace-Text "937÷5=187, 2" "939÷9=104, 3"
Replace-Text "219÷6=36, 3" "973÷9=108, 1"
Replace-Text "997÷2=498, 1" "116÷6=19, 2"

Replace-Text "845÷6=140, 5" "108÷8=13, 4"
Replace-Text "598÷8=74, 6" "765÷7=109, 2"
Replace-Text "979÷8=122, 3" "955÷2=477, 1"
Replace-Text "499÷5=99, 4" "205÷2=102, 1"
Replace-Text "149÷8=18, 5" "480÷6=80, 0"
